$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2023 column (K) to the right of the existing 2022 column (J),
# mirroring the formatting used in column J for each row.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 723.9
$ws.Range("K5").Value = 644.5
$ws.Range("K6").Value = 777.5
